$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / percentage / coin-name / link cell updates (safe as literal strings)
$ws.Range('D2').Value = '41.174.41'
$ws.Range('E2').Value = '  -1.90%  '
$ws.Range('D3').Value = '2.168.75'
$ws.Range('E3').Value = '  -1.98%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('E5').Value = '  -2.17%  '
$ws.Range('E6').Value = '  -1.62%  '
$ws.Range('E7').Value = '  -4.35%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -4.85%  '
$ws.Range('E10').Value = '  -7.48%  '
$ws.Range('E11').Value = '  -2.99%  '
$ws.Range('E12').Value = '  -4.56%  '
$ws.Range('E13').Value = '  -4.31%  '
$ws.Range('E14').Value = '  -2.30%  '
$ws.Range('D15').Value = '2.492.57'
$ws.Range('E15').Value = '  -2.08%  '
$ws.Range('E16').Value = '  -2.28%  '
$ws.Range('E17').Value = '  -4.47%  '
$ws.Range('D18').Value = '2.159.98'
$ws.Range('E18').Value = '  -2.81%  '
$ws.Range('D19').Value = '41.001.65'
$ws.Range('E19').Value = '  -1.93%  '
$ws.Range('E20').Value = '  -6.65%  '
$ws.Range('E21').Value = '  -2.82%  '
$ws.Range('E22').Value = '  -3.28%  '
$ws.Range('E23').Value = '  -6.05%  '
$ws.Range('E24').Value = '  -1.12%  '
$ws.Range('E25').Value = '  -6.63%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  -4.86%  '
$ws.Range('E28').Value = '  -1.19%  '
$ws.Range('E29').Value = '  -2.94%  '
$ws.Range('E30').Value = '  +1.41%  '
$ws.Range('E31').Value = '  +0.79%  '
$ws.Range('E32').Value = '  -2.97%  '
$ws.Range('E33').Value = '  +7.21%  '
$ws.Range('E34').Value = '  -2.90%  '
$ws.Range('E35').Value = '  -6.96%  '
$ws.Range('E36').Value = '  -2.97%  '
$ws.Range('E37').Value = '  -7.64%  '
$ws.Range('E38').Value = '  -3.03%  '
$ws.Range('E39').Value = '  -5.28%  '
$ws.Range('E40').Value = '  -6.61%  '
$ws.Range('E41').Value = '  -1.96%  '
$ws.Range('E43').Value = '  -9.15%  '
$ws.Range('E44').Value = '  -4.29%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E45').Value = '  -2.78%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E46').Value = '  -4.34%  '
$ws.Range('E47').Value = '  -5.74%  '
$ws.Range('E48').Value = '  -2.06%  '
$ws.Range('E49').Value = '  -1.98%  '
$ws.Range('E50').Value = '  -8.28%  '
$ws.Range('E51').Value = '  -2.99%  '

# Numeric-looking price text must be forced to Text format so Excel does not
# auto-convert the string into a floating point number (losing exact digits like
# trailing zeros, e.g. "236.23" -> 236.22999999999999 or "8.30" -> 8.3).
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '236.23'
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.615'
$c.Style = 'Normal'
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '70.10'
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '40.11'
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0926'
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '54.93'
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '6.77'
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '13.89'
$c.Style = 'Normal'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '70.41'
$c.Style = 'Normal'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '9.75'
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '226.26'
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '1.94'
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '10.90'
$c.Style = 'Normal'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '2.18'
$c.Style = 'Normal'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '168.20'
$c.Style = 'Normal'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '19.94'
$c.Style = 'Normal'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '30.98'
$c.Style = 'Normal'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '5.14'
$c.Style = 'Normal'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '11.98'
$c.Style = 'Normal'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '59.67'
$c.Style = 'Normal'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0976'
$c.Style = 'Normal'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '8.30'
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '97.66'
$c.Style = 'Normal'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '2.62'
$c.Style = 'Normal'
